# budget.xlsx — rework the "Composant" table:
#   - column C becomes a USD ("Coût ($)") accounting column
#   - column D is repurposed from a "Lien" (link text) column into a EUR
#     ("Coût (€)") accounting column
#   - row 3 switches from "Arduino Leonardo" to "Arduino Zero" with a new
#     price; rows 4-6 get refreshed prices; the old plain-text URL / "voir
#     OnShape" notes in column D are removed
#   - the "Total" row keeps its SUM formula and gains a (blank) euro cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usdFormat = "_-[`$`$-409]* #,##0.00_ ;_-[`$`$-409]* \-#,##0.00\ ;_-[`$`$-409]* ""-""??_ ;_-@_ "
$eurFormat = "_-* #,##0.00\ [`$€-40C]_-;\-* #,##0.00\ [`$€-40C]_-;_-* ""-""??\ [`$€-40C]_-;_-@_-"

# --- create the "Hyperlink" named cell style without disturbing any cell --
# (add + immediately drop a real hyperlink on a scratch row, then delete
# that scratch row so no trace of it remains in the used range)
$scratch = $ws.Range("B100")
$ws.Hyperlinks.Add($scratch, "https://example.com") | Out-Null
$ws.Hyperlinks.Item(1).Delete()
$scratch.EntireRow.Delete() | Out-Null

# --- header row -------------------------------------------------------------
$ws.Range("C2").Value = "Coût (`$)"
$ws.Range("D2").Value = "Coût (€)"

# --- row 3 : Arduino Leonardo -> Arduino Zero --------------------------------
$ws.Range("B3").Value = "Arduino Zero"
$ws.Range("C3").Value = 12.72
$ws.Range("D3").ClearContents()
$ws.Range("D3").Style = "Hyperlink"

# --- row 4 : Ecran LCD Tactile (price refresh, drop the link) ---------------
$ws.Range("C4").Value = 5.36
$ws.Range("D4").ClearContents()
$ws.Range("D4").ClearFormats()

# --- row 5 : Boitier haut (drop the "voir OnShape" note) --------------------
$ws.Range("D5").ClearContents()
$ws.Range("D5").ClearFormats()

# --- row 6 : Boitier bas (drop the "voir OnShape" note) ---------------------
$ws.Range("D6").ClearContents()
$ws.Range("D6").ClearFormats()

# --- apply the "Monétaire" named style to the whole C column range ----------
# (applying a named style resets any direct formatting, so fix the bold/size
# overrides back up afterwards on the header and total rows)
$ws.Range("C2").Style = "Monétaire"
$ws.Range("C3:C6").Style = "Monétaire"
$ws.Range("C9").Style = "Monétaire"

$ws.Range("C2").Font.Bold = $true
$ws.Range("C2").Font.Size = 14
$ws.Range("C9").Font.Bold = $true
$ws.Range("C9").Font.Size = 11

# --- number formats -----------------------------------------------------------
# Column C -> USD accounting format
$ws.Range("C2").NumberFormat = $usdFormat
$ws.Range("C3:C6").NumberFormat = $usdFormat
$ws.Range("C9").NumberFormat = $usdFormat

# Column D (header + D3 + total) -> EUR accounting format.
$ws.Range("D2").NumberFormat = $eurFormat
$ws.Range("D3").NumberFormat = $eurFormat
$ws.Range("D9").NumberFormat = $eurFormat

# Re-assert the bold header/total fonts that NumberFormat may have touched.
$ws.Range("D2").Font.Bold = $true
$ws.Range("D2").Font.Size = 14
$ws.Range("D9").Font.Bold = $true
$ws.Range("D9").Font.Size = 11

# Re-assert the hyperlink look on D3 (blue, underlined) after the number
# format change.
$ws.Range("D3").Font.Underline = 2
$ws.Range("D3").Font.ThemeColor = 10

# --- column widths -------------------------------------------------------------
$ws.Columns("C:D").AutoFit() | Out-Null

# --- selection -------------------------------------------------------------------
$ws.Range("D3").Select() | Out-Null
